$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 6: fill in the previously empty Category cell (A6) ---
$ws.Range("A6").Value = "Arrays"

# --- Row 7: new "Longest Substring Without Repeating Characters" entry ---
# Register the external hyperlink relationship for C7 (display text mirrors
# the URL, matching the convention already used by the other hyperlinks in
# this sheet), then set the visible cell text and restore the cell to the
# workbook's shared "Hyperlink" cell style (same style already used by the
# other hyperlink cells in the sheet).
$ws.Hyperlinks.Add($ws.Range("C7"), "https://leetcode.com/problems/longest-substring-without-repeating-characters/", [Type]::Missing, [Type]::Missing, "https://leetcode.com/problems/longest-substring-without-repeating-characters/")
$ws.Range("C7").Value = "Longest Substring Without Repeating Characters"
$ws.Range("C7").Style = "Hyperlink"

$ws.Range("A7").Value = "Arrays "
$ws.Range("B7").Value = 3
$ws.Range("D7").Value = "Brute Force - Check all the substring one by one to see if it has no duplicate character."
$ws.Range("E7").Value = "Brute Force - O(n**3)"
$ws.Range("F7").Value = "Brute Force - The space taken by the char set O(min(n,m) where n in the length of the string and m is the char set (128 for ASCII , 26 for smaller alphabets`n"

# Row height for the new row (wrapped multi-line text)
$ws.Rows.Item(7).RowHeight = 58

# Widen column C to fit the new, longer problem name
$ws.Columns.Item(3).ColumnWidth = 52

# Update the view: move the active selection
$ws.Range("D15").Select() | Out-Null
